$wb = $excel.ActiveWorkbook

# --- Add the two new worksheets at the end of the workbook --------------
$wsAdd = $wb.Worksheets.Add()
$wsAdd.Name = "Thêm sản phẩm"
$wsAdd.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$wsEdit = $wb.Worksheets.Add()
$wsEdit.Name = "Sửa sản phẩm"
$wsEdit.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch references by name now that the sheets have been moved, since
# the handles captured before Move() point at the old tab position.
$wsAdd = $wb.Worksheets.Item("Thêm sản phẩm")
$wsEdit = $wb.Worksheets.Item("Sửa sản phẩm")

# --- Populate "Thêm sản phẩm" --------------------------------------------
$wsAdd.Range("A1").Value = "username"
$wsAdd.Range("B1").Value = "password"
$wsAdd.Range("C1").Value = "name"
$wsAdd.Range("D1").Value = "image"
$wsAdd.Range("E1").Value = "title"
$wsAdd.Range("F1").Value = "brand"
$wsAdd.Range("G1").Value = "price"
$wsAdd.Range("H1").Value = "category"

$wsAdd.Range("A2").Value = "admin"
$wsAdd.Range("B2").Value = "adminpassword"
$wsAdd.Range("C2").Value = "Đồng hồ siêu cấp smartwatch"
$wsAdd.Range("D2").Value = "D:\admin.png"
$wsAdd.Range("E2").Value = "Đồng hồ"
$wsAdd.Range("F2").Value = "Casio"
$wsAdd.Range("G2").Value = 12
$wsAdd.Range("H2").Value = "Điện tử"

$wsAdd.Range("A3").Value = "admin"
$wsAdd.Range("B3").Value = "adminpassword"
$wsAdd.Range("C3").Value = "Đồng hồ siêu cấp smartwatch 2"
$wsAdd.Range("E3").Value = "Đồng hồ 2"
$wsAdd.Range("F3").Value = "Casio 2"
$wsAdd.Range("G3").Value = 13
$wsAdd.Range("H3").Value = "Điện tử"

$wsAdd.Range("A4").Value = "admin"
$wsAdd.Range("B4").Value = "adminpassword"
$wsAdd.Range("D4").Value = "D:\admin.png"

# --- Populate "Sửa sản phẩm" ----------------------------------------------
$wsEdit.Range("A1").Value = "username"
$wsEdit.Range("B1").Value = "password"
$wsEdit.Range("C1").Value = "name"
$wsEdit.Range("D1").Value = "image"
$wsEdit.Range("E1").Value = "title"
$wsEdit.Range("F1").Value = "brand"
$wsEdit.Range("G1").Value = "price"
$wsEdit.Range("H1").Value = "category"

$wsEdit.Range("A2").Value = "admin"
$wsEdit.Range("B2").Value = "adminpassword"
$wsEdit.Range("C2").Value = "Đồng hồ siêu cấp smartwatch 21212"
$wsEdit.Range("D2").Value = "D:\admin.png"
$wsEdit.Range("E2").Value = "Đồng hồ 2121"
$wsEdit.Range("F2").Value = "Casio"
$wsEdit.Range("G2").Value = 12
$wsEdit.Range("H2").Value = "Điện tử"

$wsEdit.Range("A3").Value = "admin"
$wsEdit.Range("B3").Value = "adminpassword"
$wsEdit.Range("C3").Value = "Đồng hồ siêu cấp smartwatch 2 212112"
$wsEdit.Range("E3").Value = "Đồng hồ 2"
$wsEdit.Range("F3").Value = "Casio 2"
$wsEdit.Range("G3").Value = 13
$wsEdit.Range("H3").Value = "Điện tử 21"

$wsEdit.Range("A4").Value = "admin"
$wsEdit.Range("B4").Value = "adminpassword"
$wsEdit.Range("D4").Value = "D:\admin.png"

# --- Selections / active sheet -------------------------------------------
# "Thanh toán" keeps the data it had, but loses the tab selection and gets
# a new selected cell.
$wsPay = $wb.Worksheets.Item("Thanh toán")
$wsPay.Activate()
$wsPay.Range("H12").Select()

# "Sửa sản phẩm" ends up being the active (selected) sheet/tab, with F4
# selected inside it.
$wsEdit.Activate()
$wsEdit.Range("F4").Select()
